$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph 3 is currently empty; fill it with six runs of text
# (each becomes its own <w:r> with an empty <w:rPr/>, matching the diff).
# We build the six parts separated by a CR so each lands in its own
# paragraph, then delete the paragraph marks between them so the resulting
# single paragraph keeps six distinct runs instead of Word merging them
# into one run of text.
# ---------------------------------------------------------------------------
$parts1 = @(
    'I sagen oprettet ${createdAt:d. MMMM YYYY} ',
    '(som var en ${createdAt:EEEE}) ',
    'har vi, ',
    '${board.name}',
    ', fundet at ',
    'klageren'
)
$targetIndex1 = 3
$joined1 = [string]::Join([char]13, $parts1)
$d.Paragraphs.Item($targetIndex1).Range.Text = $joined1

$n1 = $parts1.Count
for ($k = $n1 - 1; $k -ge 1; $k--) {
    $idx = $targetIndex1 + $k - 1
    $cur = $d.Paragraphs.Item($idx).Range
    $markRange = $d.Range($cur.End - 1, $cur.End)
    $markRange.Delete()
}

# ---------------------------------------------------------------------------
# Change 2: paragraph holding "${complainant}" gets a manual line break
# followed by "${complainant} ${complainant}" appended inside the SAME run.
# [char]11 is Word's internal manual-line-break character and serializes
# as <w:br/> between two <w:t> elements inside one run.
# ---------------------------------------------------------------------------
$complainantIndex = 4
$lineBreak = [char]11
$d.Paragraphs.Item($complainantIndex).Range.Text = '${complainant}' + $lineBreak + '${complainant} ${complainant}'

# ---------------------------------------------------------------------------
# Change 3: remove the now-superfluous empty paragraph that precedes the
# final URL paragraph, and give the URL paragraph explicit
# before/after spacing (0 / 283 twips = 0 / 14.15pt).
# ---------------------------------------------------------------------------
$emptyTrailingIndex = 10
$d.Paragraphs.Item($emptyTrailingIndex).Range.Delete()

$urlIndex = 10
$urlParagraph = $d.Paragraphs.Item($urlIndex)
$urlParagraph.SpaceBefore = 0
$urlParagraph.SpaceAfter = 14.15

# ---------------------------------------------------------------------------
# Change 4: the "Normal" paragraph style's default spacing after also moves
# from 0 to 283 twips (14.15pt); before stays 0.
# ---------------------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.SpaceBefore = 0
$normalStyle.ParagraphFormat.SpaceAfter = 14.15
